# Updated cryptos list values (Price and Volume(1h) columns) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.840.85"
$ws.Range("D2").Style = $ws.Range("B2").Style
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "'1.869.79"
$ws.Range("D3").Style = $ws.Range("B3").Style
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("D4").Value = "'0.9986"
$ws.Range("D4").Style = $ws.Range("B4").Style
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'0.7421"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "  -4.33%  "
$ws.Range("D6").Value = "'241.79"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("D7").Value = "'0.9986"
$ws.Range("D7").Style = $ws.Range("B7").Style
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'0.3154"
$ws.Range("D8").Style = $ws.Range("B8").Style
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("E9").Value = "  -4.30%  "
$ws.Range("D10").Value = "'0.07114"
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = "  -2.12%  "
$ws.Range("D11").Value = "'0.08380"
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = "  -6.26%  "
$ws.Range("D12").Value = "'0.7523"
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = "  -2.83%  "
$ws.Range("D13").Value = "'5.437"
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").Value = "'1.852.59"
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = "  -1.81%  "
$ws.Range("D15").Value = "'92.53"
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("D16").Value = "'29.847.01"
$ws.Range("D16").Style = $ws.Range("B16").Style
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "'6.032"
$ws.Range("D17").Style = $ws.Range("B17").Style
$ws.Range("E17").Value = "  -2.72%  "
$ws.Range("E18").Value = "  -2.98%  "
$ws.Range("D19").Value = "'242.71"
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("D20").Value = "'0.000007810"
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = "  -1.02%  "
$ws.Range("D21").Value = "'0.9988"
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "'2.116.25"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("D23").Value = "'7.932"
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("D24").Value = "'0.9992"
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'0.1572"
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D26").Value = "'9.297"
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = "  -2.68%  "
$ws.Range("D27").Value = "'164.08"
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("E28").Value = "  -1.54%  "
$ws.Range("D29").Value = "'2.021"
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("E30").Value = "  +3.14%  "
$ws.Range("E31").Value = "  +2.02%  "
$ws.Range("D32").Value = "'1.528"
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = "  -1.13%  "
$ws.Range("D33").Value = "'4.311"
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").Value = "  +4.55%  "
$ws.Range("D34").Value = "'0.05320"
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value = "  -3.73%  "
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("D36").Value = "'0.7521"
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("D37").Value = "'0.9990"
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("D39").Value = "'0.01954"
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("D40").Value = "'2.746"
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("D42").Value = "'1.103.03"
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("D43").Value = "'6.075"
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").Value = "'72.13"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = "  -2.76%  "
$ws.Range("D45").Value = "'0.8597"
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").Value = "'0.9997"
$ws.Range("D46").Style = $ws.Range("B46").Style
$ws.Range("D47").Value = "'103.18"
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").Value = "'7.679"
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").Value = "'1.841"
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Value = "  -2.98%  "
$ws.Range("D50").Value = "'3.043"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("D51").Value = "'2.015.81"
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = "  -0.04%  "
